# Applies the "repull data, push all data, mean calculation" update:
# updates the dSF column (F) values for several rows on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    8  = -2
    9  = -7
    10 = -4
    12 = -5
    13 = -2
    15 = -4
    16 = -3
    19 = -3
    23 = 3
    24 = -4
    25 = -8
    27 = 0
    31 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
